$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H138").Value = 4074.691
$ws.Range("J138").Value = 5062.6523
$ws.Range("L138").Value = 15187.9569
$ws.Range("N138").Value = -25467.9569

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5820.825
$ws.Range("I32").Value = 3273.879
$ws.Range("K32").Value = 3273.879
$ws.Range("M32").Value = -2986.879
$ws.Range("H63").Value = 3374.9167
$ws.Range("I63").Value = 2769.9
$ws.Range("K63").Value = 2769.9
$ws.Range("M63").Value = -2083.9
$ws.Range("H66").Value = 3374.9167
$ws.Range("I66").Value = 2769.9
$ws.Range("K66").Value = 13849.5
$ws.Range("M66").Value = -10417.5
$ws.Range("H122").Value = 771930.5600000001
$ws.Range("I122").Value = 835799.75
$ws.Range("K122").Value = 2507399.25
$ws.Range("M122").Value = -2504949.25
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 34056
$ws.Range("J88").Value = 34056
$ws.Range("L88").Value = 34056
$ws.Range("N88").Value = -34868
$ws.Range("H91").Value = 34056
$ws.Range("J91").Value = 34056
$ws.Range("L91").Value = 34056
$ws.Range("N91").Value = -36864
$ws.Range("H134").Value = 3820.5881
$ws.Range("I134").Value = 3390
$ws.Range("K134").Value = 10170
$ws.Range("M134").Value = -7635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3638.6086
$ws.Range("I31").Value = 3162
$ws.Range("K31").Value = 3162
$ws.Range("M31").Value = -2867
$ws.Range("H34").Value = 3638.6086
$ws.Range("I34").Value = 3162
$ws.Range("K34").Value = 3162
$ws.Range("M34").Value = -2960
$ws.Range("H68").Value = 47000
$ws.Range("J68").Value = 47000
$ws.Range("L68").Value = 47000
$ws.Range("N68").Value = -48498
$ws.Range("H71").Value = 47000
$ws.Range("J71").Value = 47000
$ws.Range("L71").Value = 141000
$ws.Range("N71").Value = -148488
$ws.Range("H74").Value = 56500
$ws.Range("J74").Value = 56500
$ws.Range("L74").Value = 56500
$ws.Range("N74").Value = -58248
$ws.Range("H77").Value = 56500
$ws.Range("J77").Value = 56500
$ws.Range("L77").Value = 169500
$ws.Range("N77").Value = -178236
$ws.Range("H86").Value = 16000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 16000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 16000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -18246
$ws.Range("H89").Value = 16000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 16000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 80000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -91232
$ws.Range("H122").Value = 3004
$ws.Range("I122").Value = 3127.2144
$ws.Range("K122").Value = 9381.643199999999
$ws.Range("M122").Value = -6931.643199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1212.5
$ws.Range("I34").Value = 283.33334
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 850.0000200000001
$ws.Range("L34").Value = 12000
$ws.Range("M34").Value = -766.0000200000001
$ws.Range("N34").Value = -12168
$ws.Range("H39").Value = 2450
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 2900
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 8700
$ws.Range("M39").Value = -5706
$ws.Range("N39").Value = -9288
$ws.Range("H55").Value = 253125
$ws.Range("I55").Value = 335833.34
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 1007500.02
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -1007323.02
$ws.Range("N55").Value = -15354
$ws.Range("H82").Value = 19243.75
$ws.Range("J82").Value = 19243.75
$ws.Range("L82").Value = 57731.25
$ws.Range("N82").Value = -58543.25
$ws.Range("H85").Value = 19243.75
$ws.Range("J85").Value = 19243.75
$ws.Range("L85").Value = 57731.25
$ws.Range("N85").Value = -60539.25
$ws.Range("H114").Value = 166.66667
$ws.Range("J114").Value = 150
$ws.Range("L114").Value = 450
$ws.Range("N114").Value = -6958
$ws.Range("H131").Value = 1095.931
$ws.Range("J131").Value = 1208.4783
$ws.Range("L131").Value = 3625.4349
$ws.Range("N131").Value = -13705.4349
$ws.Range("H140").Value = 3200
$ws.Range("I140").Value = 1400
$ws.Range("K140").Value = 4200
$ws.Range("M140").Value = 980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10062.375
$ws.Range("I20").Value = 5999
$ws.Range("J20").Value = 10642.857
$ws.Range("K20").Value = 5999
$ws.Range("L20").Value = 10642.857
$ws.Range("M20").Value = -5754
$ws.Range("N20").Value = -11132.857
$ws.Range("H70").Value = 7999.857
$ws.Range("I70").Value = 7999.6665
$ws.Range("K70").Value = 7999.6665
$ws.Range("M70").Value = -7729.6665
$ws.Range("H73").Value = 7999.857
$ws.Range("I73").Value = 7999.6665
$ws.Range("K73").Value = 7999.6665
$ws.Range("M73").Value = -7063.6665
$ws.Range("H122").Value = 38330.75
$ws.Range("I122").Value = 2240.682
$ws.Range("K122").Value = 6722.045999999999
$ws.Range("M122").Value = -4272.045999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2314.158
$ws.Range("J40").Value = 1992.75
$ws.Range("L40").Value = 1992.75
$ws.Range("N40").Value = -2264.75
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H132").Value = 4632.8335
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
$ws.Range("H136").Value = 3304.0952
$ws.Range("J136").Value = 8229
$ws.Range("L136").Value = 24687
$ws.Range("N136").Value = -29787

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7388.5
$ws.Range("J62").Value = 7599.467
$ws.Range("L62").Value = 7599.467
$ws.Range("N62").Value = -8847.467000000001
$ws.Range("H65").Value = 7388.5
$ws.Range("J65").Value = 7599.467
$ws.Range("L65").Value = 37997.335
$ws.Range("N65").Value = -44237.335
$ws.Range("H81").Value = 12475.077
$ws.Range("I81").Value = 11666
$ws.Range("K81").Value = 23332
$ws.Range("M81").Value = -22271
$ws.Range("H84").Value = 12475.077
$ws.Range("I84").Value = 11666
$ws.Range("K84").Value = 116660
$ws.Range("M84").Value = -111356
$ws.Range("H88").Value = 46666.668
$ws.Range("J88").Value = 42500
$ws.Range("L88").Value = 42500
$ws.Range("N88").Value = -43312
$ws.Range("H91").Value = 46666.668
$ws.Range("J91").Value = 42500
$ws.Range("L91").Value = 42500
$ws.Range("N91").Value = -45308
$ws.Range("H107").Value = 1027.6923
$ws.Range("I107").Value = 1030
$ws.Range("K107").Value = 3090
$ws.Range("M107").Value = -1170
